$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$v28 = $ws.Range("B28:AB28").Value()
$v29 = $ws.Range("B29:AB29").Value()
$ws.Range("B28:AB28").Value = $v29
$ws.Range("B29:AB29").Value = $v28

$v47 = $ws.Range("B47:AB47").Value()
$v48 = $ws.Range("B48:AB48").Value()
$ws.Range("B47:AB47").Value = $v48
$ws.Range("B48:AB48").Value = $v47

$v55 = $ws.Range("B55:AB55").Value()
$v56 = $ws.Range("B56:AB56").Value()
$ws.Range("B55:AB55").Value = $v56
$ws.Range("B56:AB56").Value = $v55

$v101 = $ws.Range("B101:AB101").Value()
$v102 = $ws.Range("B102:AB102").Value()
$ws.Range("B101:AB101").Value = $v102
$ws.Range("B102:AB102").Value = $v101

$v109 = $ws.Range("B109:AB109").Value()
$v110 = $ws.Range("B110:AB110").Value()
$ws.Range("B109:AB109").Value = $v110
$ws.Range("B110:AB110").Value = $v109

$v133 = $ws.Range("B133:AB133").Value()
$v134 = $ws.Range("B134:AB134").Value()
$ws.Range("B133:AB133").Value = $v134
$ws.Range("B134:AB134").Value = $v133

$v221 = $ws.Range("B221:AB221").Value()
$v222 = $ws.Range("B222:AB222").Value()
$ws.Range("B221:AB221").Value = $v222
$ws.Range("B222:AB222").Value = $v221

$v245 = $ws.Range("B245:AB245").Value()
$v246 = $ws.Range("B246:AB246").Value()
$ws.Range("B245:AB245").Value = $v246
$ws.Range("B246:AB246").Value = $v245

$v251 = $ws.Range("B251:AB251").Value()
$v252 = $ws.Range("B252:AB252").Value()
$ws.Range("B251:AB251").Value = $v252
$ws.Range("B252:AB252").Value = $v251
